$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 42
$ws1.Range("F4").Value = 244
$ws1.Range("F5").Value = 3882
$ws1.Range("F7").Value = 437

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 42
$ws4.Range("F4").Value = 244
$ws4.Range("F5").Value = 3882
$ws4.Range("F9").Value = 437

$wb.Save()
